# Applies the "Trade #102 / Trade #68 closed" update to live_trading_results.xlsx
#
# Summary of changes (per the OOXML diff):
#  - Summary sheet: OVERALL row (row 2) and momentum row (row 4) stats refreshed
#  - leadlag sheet: new OPEN trade #102 appended as row 78
#  - momentum sheet: trade #68 (row 16) transitions from OPEN -> CLOSED
#  - All Trades sheet: the now-closed trade #68 appended as a new row 69
#  - Comparison sheet: momentum row (row 3) stats refreshed
#
# NOTE: text cells that look numeric/date-like (e.g. "67.6%", "7.43", "2026-02-16")
# are written with a leading apostrophe so Excel keeps them as literal text
# instead of silently re-interpreting them as numbers/dates/percentages.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 68
$summary.Range("D2").Value = "'67.6%"
$summary.Range("E2").Value = "'+18.3220%"
$summary.Range("F2").Value = "'+0.2694%"

$summary.Range("D4").Value = "'52.0%"
$summary.Range("E4").Value = "'+7.2307%"
$summary.Range("F4").Value = "'+0.2892%"

# ---------------------------------------------------------------------------
# 2. leadlag sheet: append new OPEN trade #102 as row 78
# ---------------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Cells.Item(78, 1).Value = 102
$leadlag.Cells.Item(78, 2).Value = "'2026-02-16"
$leadlag.Cells.Item(78, 3).Value = "'21:40:21"
$leadlag.Cells.Item(78, 4).Value = "leadlag"
$leadlag.Cells.Item(78, 5).Value = "UP"
$leadlag.Cells.Item(78, 6).Value = 68412.245
# G78 (Exit Price) intentionally left blank - trade is still OPEN
$leadlag.Cells.Item(78, 8).Value = "OPEN"
$leadlag.Cells.Item(78, 9).Value = 0
$leadlag.Cells.Item(78, 10).Value = 0
$leadlag.Cells.Item(78, 11).Value = 0.75
$leadlag.Cells.Item(78, 12).Value = "Binance leading with 0.146% move"
# M78 (Exit Reason) intentionally left blank - trade is still OPEN
$leadlag.Cells.Item(78, 14).Value = 0

# ---------------------------------------------------------------------------
# 3. momentum sheet: trade #68 (row 16) OPEN -> CLOSED
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Cells.Item(16, 7).Value = 68181.761918
$momentum.Cells.Item(16, 8).Value = "CLOSED"
$momentum.Cells.Item(16, 9).Value = 0.5476
$momentum.Cells.Item(16, 10).Value = 5.48
$momentum.Cells.Item(16, 13).Value = "time_exit_5min"
$momentum.Cells.Item(16, 14).Value = 5

# ---------------------------------------------------------------------------
# 4. All Trades sheet: append the just-closed trade #68 as row 69
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(69, 1).Value = 68
$allTrades.Cells.Item(69, 2).Value = "'2026-02-16"
$allTrades.Cells.Item(69, 3).Value = "'21:35:18"
$allTrades.Cells.Item(69, 4).Value = "momentum"
$allTrades.Cells.Item(69, 5).Value = "DOWN"
$allTrades.Cells.Item(69, 6).Value = 68557.14999999999
$allTrades.Cells.Item(69, 7).Value = 68181.761918
$allTrades.Cells.Item(69, 8).Value = "CLOSED"
$allTrades.Cells.Item(69, 9).Value = 0.5476
$allTrades.Cells.Item(69, 10).Value = 5.48
$allTrades.Cells.Item(69, 11).Value = 0.9
$allTrades.Cells.Item(69, 12).Value = "Downward momentum: -0.223% over 10 samples"
$allTrades.Cells.Item(69, 13).Value = "time_exit_5min"
$allTrades.Cells.Item(69, 14).Value = 5

# ---------------------------------------------------------------------------
# 5. Comparison sheet: refresh momentum row (row 3) stats
# ---------------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("C3").Value = "'52.0%"
$comparison.Range("D3").Value = "'7.43"
$comparison.Range("E3").Value = "'+0.6427%"
$comparison.Range("G3").Value = "'1.14"
